$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Column D holds numeric-looking price strings that
# must stay text (preserve trailing zeros / thousand-dot formatting), column E
# holds percentage strings that are always non-numeric text already.
$updates = [ordered]@{
    "D2" = "23.731.11"
    "E2" = "  +1.75%  "
    "D3" = "1.652.38"
    "E3" = "  +1.61%  "
    "E4" = "  -0.24%  "
    "D5" = "0.9994"
    "E5" = "  -0.31%  "
    "D6" = "303.63"
    "E6" = "  +0.09%  "
    "D7" = "0.3820"
    "E7" = "  +2.09%  "
    "D8" = "51.24"
    "E8" = "  -0.04%  "
    "D9" = "0.3604"
    "E9" = "  -0.49%  "
    "D10" = "1.250"
    "E10" = "  +2.09%  "
    "D11" = "0.08248"
    "E11" = "  +1.25%  "
    "D12" = "0.9981"
    "E12" = "  -0.47%  "
    "D13" = "22.63"
    "E13" = "  +1.76%  "
    "E14" = "  +0.98%  "
    "E15" = "  +1.48%  "
    "E16" = "  -0.16%  "
    "D17" = "1.661.27"
    "E17" = "  +2.31%  "
    "D18" = "97.36"
    "E18" = "  +3.81%  "
    "D19" = "0.06969"
    "E19" = "  +0.34%  "
    "D20" = "6.787"
    "E20" = "  +4.02%  "
    "D21" = "17.70"
    "E21" = "  +1.26%  "
    "D22" = "1.000"
    "D23" = "12.62"
    "E23" = "  +0.74%  "
    "D24" = "23.764.28"
    "E24" = "  +1.92%  "
    "D25" = "2.544"
    "E25" = "  +3.36%  "
    "D26" = "3.079"
    "E26" = "  -1.88%  "
    "D27" = "21.33"
    "E27" = "  +0.44%  "
    "D28" = "151.55"
    "E28" = "  +0.67%  "
    "D29" = "5.265"
    "E29" = "  +0.09%  "
    "D30" = "134.86"
    "E30" = "  +1.56%  "
    "D31" = "1.841.64"
    "E31" = "  +2.43%  "
    "D32" = "6.871"
    "E32" = "  +1.86%  "
    "D33" = "1.092"
    "E33" = "  +6.05%  "
    "D34" = "11.82"
    "E34" = "  +10.57%  "
    "D35" = "2.096"
    "E35" = "  -5.91%  "
    "D36" = "0.02840"
    "E36" = "  +3.39%  "
    "D37" = "0.2519"
    "E37" = "  +1.14%  "
    "D38" = "0.08826"
    "E38" = "  +0.65%  "
    "D39" = "6.082"
    "E39" = "  +1.85%  "
    "D40" = "0.07047"
    "E40" = "  -0.83%  "
    "D41" = "12.85"
    "E41" = "  +6.20%  "
    "D42" = "0.7071"
    "E42" = "  +1.32%  "
    "D43" = "1.338"
    "E43" = "  +0.03%  "
    "E44" = "  +0.72%  "
    "D45" = "0.6535"
    "E45" = "  +0.72%  "
    "E46" = "  +3.16%  "
    "D47" = "0.9992"
    "E47" = "  -0.25%  "
    "D48" = "3.983"
    "E48" = "  +0.45%  "
    "D49" = "0.07992"
    "E49" = "  +0.37%  "
    "D50" = "128.35"
    "E50" = "  +2.56%  "
    "D51" = "1.194"
    "E51" = "  +0.52%  "
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    if ($cellRef.StartsWith("D")) {
        # Force text storage so values like "1.250" / "1.000" keep their
        # exact digits instead of being auto-coerced into numbers.
        $range.NumberFormat = "@"
        $range.Value = $updates[$cellRef]
        $range.Style = "Normal"
    } else {
        $range.Value = $updates[$cellRef]
    }
}
